$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.064374964000652
$ws.Cells.Item(2, 4).Value = 1.042326831799074
$ws.Cells.Item(2, 5).Value = 1.06883447691763
$ws.Cells.Item(2, 6).Value = 1.077623099954508
$ws.Cells.Item(2, 9).Value = 1.046091716073154
$ws.Cells.Item(2, 10).Value = 1.069335528897711
$ws.Cells.Item(2, 11).Value = 1.045103833680183
$ws.Cells.Item(2, 12).Value = 1.071538598371731
$ws.Cells.Item(2, 13).Value = 1.080303913428988
$ws.Cells.Item(2, 14).Value = 1.070854108003245

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.065827895060356
$ws.Cells.Item(3, 4).Value = 1.042950868779067
$ws.Cells.Item(3, 5).Value = 1.070131245299003
$ws.Cells.Item(3, 6).Value = 1.079008901512465
$ws.Cells.Item(3, 9).Value = 1.046430912073039
$ws.Cells.Item(3, 10).Value = 1.070441262086014
$ws.Cells.Item(3, 11).Value = 1.045539169177562
$ws.Cells.Item(3, 12).Value = 1.072650227887433
$ws.Cells.Item(3, 13).Value = 1.081506057609866
$ws.Cells.Item(3, 14).Value = 1.071961411459506

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.066767070361258
$ws.Cells.Item(4, 4).Value = 1.043354429969416
$ws.Cells.Item(4, 5).Value = 1.070969710367146
$ws.Cells.Item(4, 6).Value = 1.079905102337275
$ws.Cells.Item(4, 9).Value = 1.046648855836555
$ws.Cells.Item(4, 10).Value = 1.071155324986892
$ws.Cells.Item(4, 11).Value = 1.045819937243988
$ws.Cells.Item(4, 12).Value = 1.073368342291909
$ws.Cells.Item(4, 13).Value = 1.082282881240644
$ws.Cells.Item(4, 14).Value = 1.072676488411606

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.067161673579033
$ws.Cells.Item(5, 4).Value = 1.043524031648501
$ws.Cells.Item(5, 5).Value = 1.071322054343536
$ws.Cells.Item(5, 6).Value = 1.080281748718812
$ws.Cells.Item(5, 9).Value = 1.046740112151821
$ws.Cells.Item(5, 10).Value = 1.071455181078419
$ws.Cells.Item(5, 11).Value = 1.045937751776904
$ws.Cells.Item(5, 12).Value = 1.073669958014341
$ws.Cells.Item(5, 13).Value = 1.082609212099607
$ws.Cells.Item(5, 14).Value = 1.072976770333177

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.067227916104162
$ws.Cells.Item(6, 4).Value = 1.043552505271133
$ws.Cells.Item(6, 5).Value = 1.071381205966604
$ws.Cells.Item(6, 6).Value = 1.08034498260737
$ws.Cells.Item(6, 9).Value = 1.046755412986168
$ws.Cells.Item(6, 10).Value = 1.071505508653954
$ws.Cells.Item(6, 11).Value = 1.045957520451443
$ws.Cells.Item(6, 12).Value = 1.07372058434096
$ws.Cells.Item(6, 13).Value = 1.082663990193448
$ws.Cells.Item(6, 14).Value = 1.073027169379643

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.066772343950001
$ws.Cells.Item(7, 4).Value = 1.043356696412987
$ws.Cells.Item(7, 5).Value = 1.070974418975885
$ws.Cells.Item(7, 6).Value = 1.079910135556552
$ws.Cells.Item(7, 9).Value = 1.046650076649522
$ws.Cells.Item(7, 10).Value = 1.071159332993653
$ws.Cells.Item(7, 11).Value = 1.045821512352648
$ws.Cells.Item(7, 12).Value = 1.073372373589617
$ws.Cells.Item(7, 13).Value = 1.08228724265022
$ws.Cells.Item(7, 14).Value = 1.072680502110196

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.064866192787771
$ws.Cells.Item(8, 4).Value = 1.042537776280197
$ws.Cells.Item(8, 5).Value = 1.069272859049929
$ws.Cells.Item(8, 6).Value = 1.07809154548644
$ws.Cells.Item(8, 9).Value = 1.046206668237241
$ws.Cells.Item(8, 10).Value = 1.069709513182435
$ws.Cells.Item(8, 11).Value = 1.045251148869458
$ws.Cells.Item(8, 12).Value = 1.071914526344982
$ws.Cells.Item(8, 13).Value = 1.080710402361687
$ws.Cells.Item(8, 14).Value = 1.071228623388551

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.061499638445572
$ws.Cells.Item(9, 4).Value = 1.041092943660559
$ws.Cells.Item(9, 5).Value = 1.066269469926485
$ws.Cells.Item(9, 6).Value = 1.074882881397
$ws.Cells.Item(9, 9).Value = 1.04541349001621
$ws.Cells.Item(9, 10).Value = 1.067143682668286
$ws.Cells.Item(9, 11).Value = 1.044238996619426
$ws.Cells.Item(9, 12).Value = 1.069336366733664
$ws.Cells.Item(9, 13).Value = 1.077923625714835
$ws.Cells.Item(9, 14).Value = 1.068659149100766

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.059249749538478
$ws.Cells.Item(10, 4).Value = 1.04012850052006
$ws.Cells.Item(10, 5).Value = 1.064263562558941
$ws.Cells.Item(10, 6).Value = 1.072740742857307
$ws.Cells.Item(10, 9).Value = 1.044876673946236
$ws.Cells.Item(10, 10).Value = 1.065425431613522
$ws.Cells.Item(10, 11).Value = 1.043559413009979
$ws.Cells.Item(10, 12).Value = 1.067611127711979
$ws.Cells.Item(10, 13).Value = 1.076060010268748
$ws.Cells.Item(10, 14).Value = 1.066938457932419

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.058274132871711
$ws.Cells.Item(11, 4).Value = 1.039710591501413
$ws.Cells.Item(11, 5).Value = 1.063394056392611
$ws.Cells.Item(11, 6).Value = 1.071812389857197
$ws.Cells.Item(11, 9).Value = 1.044642305307129
$ws.Cells.Item(11, 10).Value = 1.064679525899101
$ws.Cells.Item(11, 11).Value = 1.043263993564316
$ws.Cells.Item(11, 12).Value = 1.066862490490331
$ws.Cells.Item(11, 13).Value = 1.075251617858276
$ws.Cells.Item(11, 14).Value = 1.066191492946326

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.057911528038288
$ws.Cells.Item(12, 4).Value = 1.039555316090497
$ws.Cells.Item(12, 5).Value = 1.063070937278278
$ws.Cells.Item(12, 6).Value = 1.071467433421584
$ws.Cells.Item(12, 9).Value = 1.04455495984662
$ws.Cells.Item(12, 10).Value = 1.064402173991661
$ws.Cells.Item(12, 11).Value = 1.043154087122256
$ws.Cells.Item(12, 12).Value = 1.066584168397645
$ws.Cells.Item(12, 13).Value = 1.074951124206624
$ws.Cells.Item(12, 14).Value = 1.065913747167365

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.057989317955633
$ws.Cells.Item(13, 4).Value = 1.039588625255723
$ws.Cells.Item(13, 5).Value = 1.063140254109292
$ws.Cells.Item(13, 6).Value = 1.071541433486823
$ws.Cells.Item(13, 9).Value = 1.044573708914801
$ws.Cells.Item(13, 10).Value = 1.064461680091773
$ws.Cells.Item(13, 11).Value = 1.043177670339212
$ws.Cells.Item(13, 12).Value = 1.066643880592572
$ws.Cells.Item(13, 13).Value = 1.075015591186384
$ws.Cells.Item(13, 14).Value = 1.065973337772965

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.058244164318103
$ws.Cells.Item(14, 4).Value = 1.039697757316524
$ws.Cells.Item(14, 5).Value = 1.063367350271244
$ws.Cells.Item(14, 6).Value = 1.071783878229386
$ws.Cells.Item(14, 9).Value = 1.044635091239161
$ws.Cells.Item(14, 10).Value = 1.064656605835505
$ws.Cells.Item(14, 11).Value = 1.043254912228904
$ws.Cells.Item(14, 12).Value = 1.066839489323651
$ws.Cells.Item(14, 13).Value = 1.075226783491928
$ws.Cells.Item(14, 14).Value = 1.066168540333611

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.058401154558813
$ws.Cells.Item(15, 4).Value = 1.0397649911504
$ws.Cells.Item(15, 5).Value = 1.063507252230542
$ws.Cells.Item(15, 6).Value = 1.071933239729826
$ws.Cells.Item(15, 9).Value = 1.044672872374472
$ws.Cells.Item(15, 10).Value = 1.064776667529243
$ws.Cells.Item(15, 11).Value = 1.043302480354251
$ws.Cells.Item(15, 12).Value = 1.06695997772023
$ws.Cells.Item(15, 13).Value = 1.075356876638951
$ws.Cells.Item(15, 14).Value = 1.066288772528725

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.059314467908689
$ws.Cells.Item(16, 4).Value = 1.040156229427608
$ws.Cells.Item(16, 5).Value = 1.064321248665161
$ws.Cells.Item(16, 6).Value = 1.072802337399835
$ws.Cells.Item(16, 9).Value = 1.044892187562831
$ws.Cells.Item(16, 10).Value = 1.065474894628338
$ws.Cells.Item(16, 11).Value = 1.043578994605594
$ws.Cells.Item(16, 12).Value = 1.067660778226029
$ws.Cells.Item(16, 13).Value = 1.076113629943326
$ws.Cells.Item(16, 14).Value = 1.066987991190389

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.059886985485903
$ws.Cells.Item(17, 4).Value = 1.04040156251965
$ws.Cells.Item(17, 5).Value = 1.064831593033709
$ws.Cells.Item(17, 6).Value = 1.073347283135196
$ws.Cells.Item(17, 9).Value = 1.045029242096412
$ws.Cells.Item(17, 10).Value = 1.065912363988094
$ws.Cells.Item(17, 11).Value = 1.043752134733615
$ws.Cells.Item(17, 12).Value = 1.068099940755857
$ws.Cells.Item(17, 13).Value = 1.07658793356077
$ws.Cells.Item(17, 14).Value = 1.067426081806814

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.060220790708689
$ws.Cells.Item(18, 4).Value = 1.040544632370187
$ws.Cells.Item(18, 5).Value = 1.065129178414473
$ws.Cells.Item(18, 6).Value = 1.073665064647896
$ws.Cells.Item(18, 9).Value = 1.045108998176065
$ws.Cells.Item(18, 10).Value = 1.06616735029331
$ws.Cells.Item(18, 11).Value = 1.043853013040275
$ws.Cells.Item(18, 12).Value = 1.068355943043519
$ws.Cells.Item(18, 13).Value = 1.076864448768237
$ws.Cells.Item(18, 14).Value = 1.06768143022183

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.060334587045981
$ws.Cells.Item(19, 4).Value = 1.040593410623166
$ws.Cells.Item(19, 5).Value = 1.065230632258418
$ws.Cells.Item(19, 6).Value = 1.073773407191885
$ws.Cells.Item(19, 9).Value = 1.045136161517073
$ws.Cells.Item(19, 10).Value = 1.066254263288436
$ws.Cells.Item(19, 11).Value = 1.043887391068476
$ws.Cells.Item(19, 12).Value = 1.068443207270993
$ws.Cells.Item(19, 13).Value = 1.076958710045347
$ws.Cells.Item(19, 14).Value = 1.067768466643379

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.059825573740926
$ws.Cells.Item(20, 4).Value = 1.040375243582478
$ws.Cells.Item(20, 5).Value = 1.064776847245698
$ws.Cells.Item(20, 6).Value = 1.073288823529299
$ws.Cells.Item(20, 9).Value = 1.045014556641179
$ws.Cells.Item(20, 10).Value = 1.065865446545905
$ws.Cells.Item(20, 11).Value = 1.043733569966057
$ws.Cells.Item(20, 12).Value = 1.068052838715114
$ws.Cells.Item(20, 13).Value = 1.076537059590066
$ws.Cells.Item(20, 14).Value = 1.067379097736476

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.058169124473303
$ws.Cells.Item(21, 4).Value = 1.039665621901914
$ws.Cells.Item(21, 5).Value = 1.063300480179804
$ws.Cells.Item(21, 6).Value = 1.071712487786231
$ws.Cells.Item(21, 9).Value = 1.044617023706168
$ws.Cells.Item(21, 10).Value = 1.064599213090059
$ws.Cells.Item(21, 11).Value = 1.043232171242455
$ws.Cells.Item(21, 12).Value = 1.066781894227696
$ws.Cells.Item(21, 13).Value = 1.075164598756407
$ws.Cells.Item(21, 14).Value = 1.066111066083884

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.057126387955224
$ws.Cells.Item(22, 4).Value = 1.039219191687121
$ws.Cells.Item(22, 5).Value = 1.062371383335306
$ws.Cells.Item(22, 6).Value = 1.070720658068848
$ws.Cells.Item(22, 9).Value = 1.04436539740739
$ws.Cells.Item(22, 10).Value = 1.063801405424675
$ws.Cells.Item(22, 11).Value = 1.042915911847891
$ws.Cells.Item(22, 12).Value = 1.065981381638604
$ws.Cells.Item(22, 13).Value = 1.07430039856852
$ws.Cells.Item(22, 14).Value = 1.065312125440105

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.05767928421304
$ws.Cells.Item(23, 4).Value = 1.039455877883784
$ws.Cells.Item(23, 5).Value = 1.062863997156344
$ws.Cells.Item(23, 6).Value = 1.071246516380896
$ws.Cells.Item(23, 9).Value = 1.04449894916323
$ws.Cells.Item(23, 10).Value = 1.064224498984128
$ws.Cells.Item(23, 11).Value = 1.043083663001788
$ws.Cells.Item(23, 12).Value = 1.066405884693693
$ws.Cells.Item(23, 13).Value = 1.074758650467026
$ws.Cells.Item(23, 14).Value = 1.065735819840942

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.059853323481271
$ws.Cells.Item(24, 4).Value = 1.040387136066594
$ws.Cells.Item(24, 5).Value = 1.064801584788809
$ws.Cells.Item(24, 6).Value = 1.073315239144417
$ws.Cells.Item(24, 9).Value = 1.045021192939637
$ws.Cells.Item(24, 10).Value = 1.065886647081698
$ws.Cells.Item(24, 11).Value = 1.043741958929846
$ws.Cells.Item(24, 12).Value = 1.068074122575349
$ws.Cells.Item(24, 13).Value = 1.076560047774547
$ws.Cells.Item(24, 14).Value = 1.067400328379461

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.062370919694171
$ws.Cells.Item(25, 4).Value = 1.041466681467129
$ws.Cells.Item(25, 5).Value = 1.067046541040419
$ws.Cells.Item(25, 6).Value = 1.07571291062039
$ws.Cells.Item(25, 9).Value = 1.045619955675213
$ws.Cells.Item(25, 10).Value = 1.067808348964751
$ws.Cells.Item(25, 11).Value = 1.044501507829701
$ws.Cells.Item(25, 12).Value = 1.070004004105082
$ws.Cells.Item(25, 13).Value = 1.07864507091942
$ws.Cells.Item(25, 14).Value = 1.069324759299611
